# Rounds every numeric cell in row 5 (columns B:AH) to 2 decimal places
# using "round half to even" on the *decimal* representation of the
# original value (this matches how the source dataset was produced,
# and is NOT the same as .NET's Math.Round on the binary double), and
# removes the last data row (row 6), shrinking the used range.

function RoundHalfEvenStr($s, $nd) {
    # Ensure there is a fractional part to work with.
    if ($s.IndexOf(".") -lt 0) {
        $s = "{0}.00" -f $s
    }

    $parts = $s.Split(".")
    $intpart = $parts[0]
    $fracpart = $parts[1]

    # Pad the fractional part if it is shorter than the digits we keep.
    if ($fracpart.Length -lt $nd) {
        $padCount = $nd - $fracpart.Length
        $zeros = "0" * $padCount
        $fracpart = "{0}{1}" -f $fracpart, $zeros
    }

    $keep = $fracpart.Substring(0, $nd)
    $rest = $fracpart.Substring($nd)
    if ($rest.Length -eq 0) {
        $rest = "0"
    }

    $firstRestDigit = $rest.Substring(0, 1)
    $roundUp = $false

    if ($firstRestDigit -gt "5") {
        $roundUp = $true
    } elseif ($firstRestDigit -lt "5") {
        $roundUp = $false
    } else {
        # Exactly "5...": round up only if something nonzero follows,
        # otherwise it's a true tie -> round to even.
        $remainderNonzero = $false
        if ($rest.Length -gt 1) {
            $tail = $rest.Substring(1)
            if ([long]$tail -ne 0) {
                $remainderNonzero = $true
            }
        }
        if ($remainderNonzero) {
            $roundUp = $true
        } else {
            $lastKeepDigit = $keep.Substring($keep.Length - 1, 1)
            $roundUp = (([int]$lastKeepDigit) % 2) -eq 1
        }
    }

    $combinedStr = "{0}{1}" -f $intpart, $keep
    $combined = [long]$combinedStr
    if ($roundUp) {
        $combined = $combined + 1
    }

    $divisor = [Math]::Pow(10, $nd)
    $val = $combined / $divisor
    return $val
}

function GetCellValueStr($cell) {
    # Prefix with a non-numeric character so PowerShell keeps the
    # shortest round-trip decimal text instead of re-coercing it back
    # to a (possibly noisier) number.
    $v = $cell.Value2
    $s = "v$v"
    return $s.Substring(1)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5, columns B (2) through AH (34): round every value to 2 decimals.
for ($c = 2; $c -le 34; $c++) {
    $cell = $ws.Cells.Item(5, $c)
    $s = GetCellValueStr $cell
    $rounded = RoundHalfEvenStr $s 2
    $cell.Value = $rounded
}

# Drop the last data row entirely (row 6); this also shrinks the sheet
# dimension from A1:AH6 down to A1:AH5 automatically.
$ws.Rows.Item(6).Delete()
